$wb = $excel.ActiveWorkbook

# ALC row 76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 5663.026
$ws.Range("I76").Value = 4877.5
$ws.Range("J76").Value = 6336.3335
$ws.Range("K76").Value = 4877.5
$ws.Range("L76").Value = 6336.3335
$ws.Range("M76").Value = -4562.5
$ws.Range("N76").Value = -6966.3335

# ALC row 79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 5663.026
$ws.Range("I79").Value = 4877.5
$ws.Range("J79").Value = 6336.3335
$ws.Range("K79").Value = 4877.5
$ws.Range("L79").Value = 6336.3335
$ws.Range("M79").Value = -3785.5
$ws.Range("N79").Value = -8520.333500000001

# ALC row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 2495.1177
$ws.Range("I100").Value = 1868.5555
$ws.Range("J100").Value = 3200
$ws.Range("K100").Value = 1868.5555
$ws.Range("L100").Value = 3200
$ws.Range("M100").Value = -1327.5555
$ws.Range("N100").Value = -4282

# ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1199.4054
$ws.Range("I112").Value = 450
$ws.Range("J112").Value = 1290.2424
$ws.Range("K112").Value = 1350
$ws.Range("L112").Value = 3870.7272
$ws.Range("M112").Value = -242
$ws.Range("N112").Value = -6086.7272

# ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 2192.65
$ws.Range("I129").Value = 653.25
$ws.Range("J129").Value = 2577.5
$ws.Range("K129").Value = 1959.75
$ws.Range("L129").Value = 7732.5
$ws.Range("M129").Value = 3040.25
$ws.Range("N129").Value = -17732.5

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2008.25
$ws.Range("I138").Value = 1351.2916
$ws.Range("J138").Value = 2993.6875
$ws.Range("K138").Value = 4053.8748
$ws.Range("L138").Value = 8981.0625
$ws.Range("M138").Value = 1086.1252
$ws.Range("N138").Value = -19261.0625

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9857.959999999999
$ws.Range("I32").Value = 10223.777
$ws.Range("K32").Value = 10223.777
$ws.Range("M32").Value = -9936.777

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 14709094
$ws.Range("I74").Value = 23811478
$ws.Range("J74").Value = 5243.5386
$ws.Range("K74").Value = 23811478
$ws.Range("L74").Value = 5243.5386
$ws.Range("M74").Value = -23810604
$ws.Range("N74").Value = -6991.5386

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 14709094
$ws.Range("I77").Value = 23811478
$ws.Range("J77").Value = 5243.5386
$ws.Range("K77").Value = 119057390
$ws.Range("L77").Value = 26217.693
$ws.Range("M77").Value = -119053022
$ws.Range("N77").Value = -34953.693

# ARM row 109
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H109").Value = 29583.666
$ws.Range("J109").Value = 29583.666
$ws.Range("L109").Value = 29583.666
$ws.Range("N109").Value = -32357.666

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1377.875
$ws.Range("I110").Value = 1500
$ws.Range("K110").Value = 1500
$ws.Range("M110").Value = 545

# ARM row 112
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H112").Value = 13380
$ws.Range("J112").Value = 13380
$ws.Range("L112").Value = 13380
$ws.Range("N112").Value = -16334

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 6946590.5
$ws.Range("I132").Value = 11906624
$ws.Range("J132").Value = 2543.4666
$ws.Range("K132").Value = 35719872
$ws.Range("L132").Value = 7630.399800000001
$ws.Range("M132").Value = -35717342
$ws.Range("N132").Value = -12690.3998

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1333.3334
$ws.Range("J99").Value = 1500
$ws.Range("L99").Value = 1500
$ws.Range("N99").Value = -4496

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2498.25
$ws.Range("I107").Value = 2693
$ws.Range("J107").Value = 2108.75
$ws.Range("K107").Value = 2693
$ws.Range("L107").Value = 2108.75
$ws.Range("M107").Value = -773
$ws.Range("N107").Value = -5948.75

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1461.1111
$ws.Range("I99").Value = 1440
$ws.Range("J99").Value = 1487.5
$ws.Range("K99").Value = 1440
$ws.Range("L99").Value = 1487.5
$ws.Range("M99").Value = 58
$ws.Range("N99").Value = -4483.5

# CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1622.3077
$ws.Range("I122").Value = 1603.6666
$ws.Range("J122").Value = 1664.25
$ws.Range("K122").Value = 4810.9998
$ws.Range("L122").Value = 4992.75
$ws.Range("M122").Value = -2360.9998
$ws.Range("N122").Value = -9892.75

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 1461.1111
$ws.Range("I126").Value = 1440
$ws.Range("J126").Value = 1487.5
$ws.Range("K126").Value = 4320
$ws.Range("L126").Value = 4462.5
$ws.Range("M126").Value = -1850
$ws.Range("N126").Value = -9402.5

# CUL row 64
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 3381.4
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 3381.4
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 10144.2
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -10684.2

# CUL row 67
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H67").Value = 3381.4
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 3381.4
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 10144.2
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -12016.2

# CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 894.475
$ws.Range("I122").Value = 815.3333
$ws.Range("J122").Value = 981.9474
$ws.Range("K122").Value = 7337.9997
$ws.Range("L122").Value = 8837.526600000001
$ws.Range("M122").Value = -4887.9997
$ws.Range("N122").Value = -13737.5266

# CUL row 129
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 3337.348
$ws.Range("I129").Value = 1940
$ws.Range("J129").Value = 3631.5264
$ws.Range("K129").Value = 5820
$ws.Range("L129").Value = 10894.5792
$ws.Range("M129").Value = -820
$ws.Range("N129").Value = -20894.5792

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 881.7083
$ws.Range("J131").Value = 881.55
$ws.Range("L131").Value = 2644.65
$ws.Range("N131").Value = -12724.65

# CUL row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1084.1666
$ws.Range("I132").Value = 590
$ws.Range("J132").Value = 1776
$ws.Range("K132").Value = 5310
$ws.Range("L132").Value = 15984
$ws.Range("M132").Value = -2780
$ws.Range("N132").Value = -21044

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4529.385
$ws.Range("I132").Value = 3720.7827
$ws.Range("J132").Value = 5691.75
$ws.Range("K132").Value = 11162.3481
$ws.Range("L132").Value = 17075.25
$ws.Range("M132").Value = -8632.348100000001
$ws.Range("N132").Value = -22135.25

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 35727000
$ws.Range("I136").Value = 83335520
$ws.Range("J136").Value = 20613.5
$ws.Range("K136").Value = 250006560
$ws.Range("L136").Value = 61840.5
$ws.Range("M136").Value = -250004010
$ws.Range("N136").Value = -66940.5

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2397.7083
$ws.Range("I126").Value = 1736.8422
$ws.Range("J126").Value = 4909
$ws.Range("K126").Value = 5210.5266
$ws.Range("L126").Value = 14727
$ws.Range("M126").Value = -2740.5266
$ws.Range("N126").Value = -19667
